$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete rows 16-18 (data now only spans to row 15)
$ws.Range("A16:H18").Delete()

# Column A holds numeric-looking IDs that must stay text, like the rest of the sheet
$ws.Range("A2:A15").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "1329569"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1329569"
$ws.Range("C2").Value = "Accelerate Romania - Map Design Assistant for Unity Casual Game"
$ws.Range("D2").Value = "Cluj-Napoca, Romania"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "1 applicant"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Trafiki"

# Row 3
$ws.Range("A3").Value = "1329557"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1329557"
$ws.Range("C3").Value = "Taste Hungary| Associate Technical Support - Swedish"
$ws.Range("D3").Value = "Budapeste, Hungria"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "4 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Tech Mahindra Kft."

# Row 4
$ws.Range("A4").Value = "1329556"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1329556"
$ws.Range("C4").Value = "Taste Hungary| Associate Technical Support - Dutch"
$ws.Range("D4").Value = "Budapeste, Hungria"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "4 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Tech Mahindra Kft."

# Row 5
$ws.Range("A5").Value = "1329555"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1329555"
$ws.Range("C5").Value = "Taste Hungary| Associate Technical Support - Czheco"
$ws.Range("D5").Value = "Budapeste, Hungria"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "2 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "Tech Mahindra Kft."

# Row 6
$ws.Range("A6").Value = "1328965"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1328965"
$ws.Range("C6").Value = "Account Manager (German Only)"
$ws.Range("D6").Value = "Assen, Nederland"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "ICT Specialist"

# Row 7
$ws.Range("A7").Value = "1328629"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1328629"
$ws.Range("C7").Value = "Marketing Intern"
$ws.Range("D7").Value = "Hyderabad, Telangana, India"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "0 applicants"
$ws.Range("G7").Value = "9 - 12 Weeks"
$ws.Range("H7").Value = "Amaavi Luxe Travels"

# Row 8
$ws.Range("A8").Value = "1328155"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1328155"
$ws.Range("C8").Value = "Clinical Study Analyst Trainee"
$ws.Range("D8").Value = "Bruxelles, Belgio"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "52 applicants"
$ws.Range("G8").Value = "6 - 18 Months"
$ws.Range("H8").Value = "UCB"

# Row 9
$ws.Range("A9").Value = "1327813"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1327813"
$ws.Range("C9").Value = "Nursery Spanish Practitioner"
$ws.Range("D9").Value = "Ashby-de-la-Zouch LE65, UK"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "35 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "Bilingual Day Nursery and Preschool Ltd"

# Row 10
$ws.Range("A10").Value = "1327381"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1327381"
$ws.Range("C10").Value = "Product Management Intern"
$ws.Range("D10").Value = "Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "88 applicants"
$ws.Range("G10").Value = "6 - 18 Months"
$ws.Range("H10").Value = "ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ"

# Row 11
$ws.Range("A11").Value = "1327380"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1327380"
$ws.Range("C11").Value = "Comunication Intern"
$ws.Range("D11").Value = "Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "84 applicants"
$ws.Range("G11").Value = "6 - 18 Months"
$ws.Range("H11").Value = "ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ"

# Row 12
$ws.Range("A12").Value = "1326670"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1326670"
$ws.Range("C12").Value = "TIM Operations Assistant Intern"
$ws.Range("D12").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "116 applicants"
$ws.Range("G12").Value = "6 - 18 Months"
$ws.Range("H12").Value = "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"

# Row 13
$ws.Range("A13").Value = "1317292"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1317292"
$ws.Range("C13").Value = "[Impact Florianópolis]- Social Media"
$ws.Range("D13").Value = "São Miguel do Oeste, SC, 89900-000, Brasil"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "86 applicants"
$ws.Range("G13").Value = "9 - 12 Weeks"
$ws.Range("H13").Value = "KNN Idiomas"

# Row 14
$ws.Range("A14").Value = "1303804"
$ws.Range("B14").Value = "https://aiesec.org/opportunity/global-talent/1303804"
$ws.Range("C14").Value = "[Impact Florianópolis]- Social Media"
$ws.Range("D14").Value = "São Miguel do Oeste, SC, 89900-000, Brasil"
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "73 applicants"
$ws.Range("G14").Value = "9 - 12 Weeks"
$ws.Range("H14").Value = "KNN Idiomas"

# Row 15
$ws.Range("A15").Value = "1289378"
$ws.Range("B15").Value = "https://aiesec.org/opportunity/global-talent/1289378"
$ws.Range("C15").Value = "Medical Advisor (Spanish Speaker)"
$ws.Range("D15").Value = "İstanbul, Türkiye"
$ws.Range("E15").Value = "No"
$ws.Range("F15").Value = "120 applicants"
$ws.Range("G15").Value = "6 - 18 Months"
$ws.Range("H15").Value = "International Plus"

# Column width adjustments (ColumnWidth offsets by ~0.8333 on save, so compensate)
$ws.Columns.Item(3).ColumnWidth = 66 - 5/6
$ws.Columns.Item(4).ColumnWidth = 57 - 5/6
$ws.Columns.Item(8).ColumnWidth = 61 - 5/6

